$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Roraima
$ws.Range("A2").Value = "Roraima"
$ws.Range("B2").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C2").Value = 4.538560136722921

# Row 3: Mato Grosso
$ws.Range("A3").Value = "Mato Grosso"
$ws.Range("B3").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C3").Value = 3.533948479222659

# Row 4: Rio Grande do Norte
$ws.Range("A4").Value = "Rio Grande do Norte"
$ws.Range("B4").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C4").Value = 2.579552159166688

# Row 5: Acre
$ws.Range("A5").Value = "Acre"
$ws.Range("B5").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C5").Value = 2.45357712462976

# Row 6: Sergipe
$ws.Range("A6").Value = "Sergipe"
$ws.Range("B6").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C6").Value = 2.257560303510999

# Row 7: Paraíba
$ws.Range("A7").Value = "Paraíba"
$ws.Range("B7").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C7").Value = 1.909986483286929

# Row 8: Nordeste (name unchanged)
$ws.Range("B8").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C8").Value = 0.4560172374823068

# Row 9: Brasil (name unchanged)
$ws.Range("B9").Value = "Diferença 2024/02 - 2023/02"
$ws.Range("C9").Value = 0.5107653626250297
